$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 2

$ws.Range("B4").Select()
